# The workbook gains one new daily price record for "Haba" (Feria Lagunitas
# de Puerto Montt). It is inserted as a brand-new row 24, which pushes every
# existing record that used to live at rows 24-127 down by one row (to
# 25-128), exactly like Excel's own "Insert Row" behaviour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 24..127 down to 25..128, opening up a blank row 24.
$ws.Rows(24).Insert()

# Populate the newly opened row 24 with the new record.
$ws.Range("A24").Value = 4
$ws.Range("B24").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C24").Value = "Los Lagos"
$ws.Range("D24").Value = 45069
$ws.Range("E24").Value = 10
$ws.Range("F24").Value = 100112026
$ws.Range("G24").Value = "Haba"
$ws.Range("H24").Value = "Sin especificar"
$ws.Range("I24").Value = "Primera"
$ws.Range("J24").Value = 80
$ws.Range("K24").Value = 23000
$ws.Range("L24").Value = 23000
$ws.Range("M24").Value = 23000
$ws.Range("N24").Value = "$/saco 25 kilos"
$ws.Range("O24").Value = "Provincia de Limarí"
$ws.Range("P24").Value = 920
$ws.Range("Q24").Value = 25
$ws.Range("R24").Value = "Hortaliza"
